$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Turn off iterative calculation (calcPr iterate attribute removed in target)
$excel.Iteration = $false

# Clear the date stamp that used to live in C1 (remove value and formatting)
$ws.Range("C1").Clear()

# Update text referencing the currency-year baseline from 2019 -> 2020
$ws.Range("A24").Value = "2020 dollars"
$ws.Range("A21").Value = "million 2020 dollars"
$ws.Range("A18").Value = "billion 2020 dollars"
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2020 dollar."'
$ws.Range("B30").Value = "2012 dollars are worth more than 2020 dollars, so we need a"

# Updated conversion factor value
$ws.Range("A26").Value = 0.88711067149387013

# Move the active selection to match the saved view state
$ws.Range("B31").Select()
